$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 341; this pushes the existing rows
# 341-345 down to become rows 343-347 (format of row above is inherited,
# matching the date-style "s=2" used in column D).
$ws.Range("A341:A342").EntireRow.Insert()

# --- New row 341: Coliflor, Primera, week of 2021-09-09 (serial 44448) ---
$ws.Cells.Item(341, 1).Value = 3
$ws.Cells.Item(341, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(341, 3).Value = "Coquimbo"
$ws.Cells.Item(341, 4).Value = 44448
$ws.Cells.Item(341, 5).Value = 5
$ws.Cells.Item(341, 6).Value = 100112008
$ws.Cells.Item(341, 7).Value = "Coliflor"
$ws.Cells.Item(341, 8).Value = "Sin especificar"
$ws.Cells.Item(341, 9).Value = "Primera"
$ws.Cells.Item(341, 10).Value = 2950
$ws.Cells.Item(341, 11).Value = 550
$ws.Cells.Item(341, 12).Value = 600
$ws.Cells.Item(341, 13).Value = 573
$ws.Cells.Item(341, 14).Value = "$/unidad"
$ws.Cells.Item(341, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(341, 16).Value = 573
$ws.Cells.Item(341, 17).Value = 1
$ws.Cells.Item(341, 18).Value = "Hortaliza"

# --- New row 342: Coliflor, Segunda, week of 2021-09-09 (serial 44448) ---
$ws.Cells.Item(342, 1).Value = 3
$ws.Cells.Item(342, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(342, 3).Value = "Coquimbo"
$ws.Cells.Item(342, 4).Value = 44448
$ws.Cells.Item(342, 5).Value = 5
$ws.Cells.Item(342, 6).Value = 100112008
$ws.Cells.Item(342, 7).Value = "Coliflor"
$ws.Cells.Item(342, 8).Value = "Sin especificar"
$ws.Cells.Item(342, 9).Value = "Segunda"
$ws.Cells.Item(342, 10).Value = 1200
$ws.Cells.Item(342, 11).Value = 500
$ws.Cells.Item(342, 12).Value = 500
$ws.Cells.Item(342, 13).Value = 500
$ws.Cells.Item(342, 14).Value = "$/unidad"
$ws.Cells.Item(342, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(342, 16).Value = 500
$ws.Cells.Item(342, 17).Value = 1
$ws.Cells.Item(342, 18).Value = "Hortaliza"

$ws.Range("A1").Select()
